# Automatische test-sync: 2025-08-28 21:29:50
#
# Adds a new "Logs" row (row 32) for a new mail log entry, and re-syncs the
# "Dashboard" category counts to match (shifting the category labels for
# rows 4-8 and bumping the "Documentatie / Datasheets" count from 1 to 2).

$wb = $excel.ActiveWorkbook

# ---- 1. Logs sheet: append new log entry on row 32 ----
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Cells.Item(32, 1).Value = "Opvolging contact"
$wsLogs.Cells.Item(32, 2).Value = "mailmind.test@zohomail.eu"
$wsLogs.Cells.Item(32, 4).Value = "Overig"
$wsLogs.Cells.Item(32, 6).Value = "2025-08-28 21:28:55"
$wsLogs.Cells.Item(32, 7).Value = "Nee"
$wsLogs.Cells.Item(32, 8).Value = "Ja"
$wsLogs.Cells.Item(32, 9).Value = "Nee"
$wsLogs.Cells.Item(32, 10).Value = "Nee"

# Expand the conditional-formatting ranges (D, G, H, I, J) so the new row
# (32) is covered too, mirroring what Excel does when you extend a table.
$cfCols = @("D", "G", "H", "I", "J")
foreach ($col in $cfCols) {
    $oldRange = $wsLogs.Range($col + "2:" + $col + "31")
    $newRange = $wsLogs.Range($col + "2:" + $col + "32")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---- 2. Dashboard sheet: re-sync category counts/labels ----
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Cells.Item(4, 1).Value = "Overig"
$wsDash.Cells.Item(5, 1).Value = "Planning / Afspraak"
$wsDash.Cells.Item(6, 1).Value = "Kwaliteit / Certificaten"
$wsDash.Cells.Item(8, 1).Value = "Documentatie / Datasheets"
$wsDash.Cells.Item(8, 2).Value = 2
